$d = $word.ActiveDocument

# 1. Locate "sive on tablet, desktop, mobile" and compute the split point
#    right before ", mobile".
$findRng = $d.Content
$findRng.Find.Execute("sive on tablet, desktop, mobile") | Out-Null
$splitPoint = $findRng.Start + 23

# 2. Move the "_GoBack" bookmark to that split point. Word only allows a
#    single bookmark per name, so re-adding "_GoBack" here removes it from
#    its old location (end of the "Education..." paragraph) and places it
#    between "desktop" and ", mobile" -- which also splits the run in two,
#    exactly as in the target edit.
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Safety net: ensure no stray "_GoBack" bookmark remains anywhere else.
if ($d.Bookmarks.Exists("_GoBack") -and $d.Bookmarks("_GoBack").Start -ne $splitPoint) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3. Update the page margins (top 540 -> 450 twips, bottom 540 -> 720 twips).
#    Word's PageSetup margin properties are expressed in points (1 pt = 20
#    twips), so divide the target twip values by 20.
$d.Sections(1).PageSetup.TopMargin = 450 / 20
$d.Sections(1).PageSetup.BottomMargin = 720 / 20
